$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full A1:D26 table (header row untouched). Column A = index, B = Massa, C = Long, D = Lat.
$data = @(
    @(0,  2.600125879871367, -48.71384122723328, -17.51148865412093),
    @(1,  2.128184125780414, -48.6960545846077,  -17.76050165087074),
    @(2,  0.8842618317366855,-48.68716126329491, -17.75160832955824),
    @(3,  2.112953020975302, -48.65158797804375, -17.44034208362099),
    @(4,  1.999954430009206, -48.67826794198212, -17.57374190330838),
    @(5,  1.227290394705522, -48.67826794198212, -17.6804617590583),
    @(6,  4.336656630291031, -48.68716126329491, -17.75160832955824),
    @(7,  2.429556842686575, -48.64269465673096, -17.75160832955824),
    @(8,  1.904204420047018, -48.65158797804375, -17.75160832955824),
    @(9,  3.104010962916469, -48.66048129935654, -17.78718161480822),
    @(10, 1.14687466333332,  -48.67826794198212, -17.75160832955824),
    @(11, 3.349726232261484, -48.65158797804375, -17.75160832955824),
    @(12, 2.875132332161531, -48.66937462066933, -17.74271500824575),
    @(13, 0.5625408548920923,-48.72273454854607, -17.78718161480822),
    @(14, 2.078663903585458, -48.6960545846077,  -17.57374190330838),
    @(15, 2.065614962453758, -48.71384122723328, -17.76050165087074),
    @(16, 3.550036750560322, -48.66937462066933, -17.76050165087074),
    @(17, 1.672601942801667, -48.72273454854607, -17.79607493612071),
    @(18, 2.383388957533608, -48.68716126329491, -17.79607493612071),
    @(19, 2.169717741959055, -48.71384122723328, -17.75160832955824),
    @(20, 2.038109824454579, -48.67826794198212, -17.77828829349572),
    @(21, 2.556291383413752, -48.66048129935654, -17.75160832955824),
    @(22, 3.631499646315595, -48.68716126329491, -17.78718161480822),
    @(23, 1.875028941637016, -48.70494790592049, -17.76050165087074),
    @(24, 1.982599194667287, -48.73162786985886, -17.78718161480822)
)

$rowIndex = 2
foreach ($rec in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $rec[0]
    $ws.Cells.Item($rowIndex, 2).Value = $rec[1]
    $ws.Cells.Item($rowIndex, 3).Value = $rec[2]
    $ws.Cells.Item($rowIndex, 4).Value = $rec[3]

    $rowIndex++
}

# The newly appended rows (A12:A26) need the same formatting as the
# existing index column (A2:A11) - copy formats only.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A12:A26").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false
